# Fixed Bento 80 Test scripts
# Append an "ORDER BY ... LIMIT 100" clause to each of the three Cypher
# queries stored on the "startup" sheet (CasesTab / SamplesTab / FilesTab
# rows), then move the active selection to B2 (the Cases query cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - CasesTab query (column B)
$casesQuery = $ws.Range("B2").Value()
$ws.Range("B2").Value = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"

# Row 3 - SamplesTab query (column B)
$samplesQuery = $ws.Range("B3").Value()
$ws.Range("B3").Value = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"

# Row 4 - FilesTab query (column B)
$filesQuery = $ws.Range("B4").Value()
$ws.Range("B4").Value = $filesQuery + "`n order By f.file_name ASC LIMIT 100"

# Move the active selection from D2 to B2, matching the saved view state.
$ws.Range("B2").Select()
